# Case_4_120 res_line pl_mw: refreshed per-row line-flow results for the 380 kV case.
# Columns A,F,H,K are unchanged (A is the row index, F/H/K remain 0).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One array entry per data row: row number followed by the new values for
# columns B, C, D, E, G, I, J, L, M, N, O (column indices 2,3,4,5,7,9,10,12,13,14,15).
$newData = @(
    ,(2, 1.501692279645511, 0.1124359296681519, 0.4602920812456688, 0.1374608553934742, 0.002515614509111651, 1.245871032244658, 0.05277358531352405, 0.4573435427471964, 0.4109732708810228, 1.839380720024081, 5.663734183809282)
    ,(3, 1.423870602517582, 0.09873227098901793, 0.4604055823394049, 0.1382565314790041, 0.002518857635303794, 1.254432228943763, 0.05188228727547894, 0.4541975564808354, 0.3979044180125655, 1.85773601711918, 5.678321474904408)
    ,(4, 1.376568396255436, 0.09027729789096384, 0.4606569422611955, 0.1387783211297791, 0.002520957021512313, 1.260275576457769, 0.05132687694776195, 0.4524421558532481, 0.3900332372893729, 1.869622494322012, 5.690166051093229)
    ,(5, 1.357414485448572, 0.08682165277384968, 0.4608051330775425, 0.1389993278085297, 0.002521839801745741, 1.262804328537559, 0.05109849519933718, 0.4517712541051893, 0.3868644180382148, 1.874621350817119, 5.695719031355821)
    ,(6, 1.354241406415468, 0.08624723402621726, 0.4608325060948317, 0.1390365319400955, 0.002521988035915427, 1.263233137016528, 0.05106044897116035, 0.4516625388640136, 0.3863405855652999, 1.875460771803013, 5.696684962341436)
    ,(7, 1.376309583767323, 0.09023073491309219, 0.4606587554329877, 0.1387812677813294, 0.00252096881655639, 1.260309082708137, 0.0513238051997007, 0.452432927750209, 0.3899903442650512, 1.869689282898449, 5.690238000064653)
    ,(8, 1.474760345682682, 0.1077194484501263, 0.4602935561894128, 0.1377283139260443, 0.00251671035379488, 1.248701154871959, 0.05246795807444826, 0.4562222960473861, 0.4064354838739348, 1.845581760175705, 5.668164498895976)
    ,(9, 1.671590539382919, 0.1416875575050085, 0.4610159011814545, 0.1359265810849237, 0.002509213436662387, 1.230593959584411, 0.05464698556870218, 0.4650472539212558, 0.4398909193712299, 1.803193495740185, 5.647797076676113)
    ,(10, 1.81845450742037, 0.1664428960120574, 0.4624201946827355, 0.1347623468544104, 0.00250422075908031, 1.220130175184977, 0.0562086517862781, 0.4723760342987049, 0.465197550793512, 1.775024727597479, 5.646818951898439)
    ,(11, 1.885747469606599, 0.1776609021820832, 0.4632480702198478, 0.1342671503615547, 0.002502060236388104, 1.215986848062244, 0.05691060728812403, 0.4758926662656648, 0.476866396477611, 1.762854751323491, 5.649414522001308)
    ,(12, 1.911298077970287, 0.1819025685227587, 0.4635886849880961, 0.1340845674813398, 0.002501257932054029, 1.214506584844656, 0.0571752037827693, 0.4772504994427891, 0.4813074220961155, 1.758338883314298, 5.650834795933974)
    ,(13, 1.90579228442175, 0.1809893340652593, 0.4635141226069948, 0.1341236705775462, 0.002501430019450109, 1.214821439242236, 0.05711827245050216, 0.4769569040955588, 0.4803499805354932, 1.759307335959484, 5.650509458934096)
    ,(14, 1.887848173924624, 0.178009994556021, 0.4632755497492553, 0.1342520302771169, 0.002501993913179142, 1.21586328748343, 0.05693240024553248, 0.4760038523749586, 0.4772313170977753, 1.762481371460385, 5.649522601936212)
    ,(15, 1.876865732993849, 0.1761842326882572, 0.4631329462623768, 0.1343312969190844, 0.002502341374766677, 1.216513005528171, 0.05681838933485039, 0.47542348417862, 0.4753239406848806, 1.764437624463561, 5.648975089272028)
    ,(16, 1.814066374297511, 0.1657088909655045, 0.4623698893852861, 0.1347954006732115, 0.002504364174827534, 1.22041336469664, 0.05616260695645536, 0.4721498812742766, 0.4644380966736463, 1.775833026757869, 5.646710522363236)
    ,(17, 1.775663956850678, 0.1592714273241427, 0.4619501506631991, 0.1350889202435672, 0.002505633388857478, 1.222964091641067, 0.05575813816008335, 0.4701883516796244, 0.4577999480845634, 1.782988722933684, 5.646100150831387)
    ,(18, 1.753621520461309, 0.1555646980189351, 0.4617265269874906, 0.1352609853866236, 0.002506373828357498, 1.224489246425286, 0.05552470283477007, 0.4690773383278071, 0.4539966279082108, 1.787165118394839, 5.646035259700739)
    ,(19, 1.746166206349812, 0.1543089677015814, 0.4616538709651223, 0.1353198006111, 0.002506626320826747, 1.22501560543877, 0.05544552909874767, 0.4687041270818355, 0.4527114330048221, 1.788589586827193, 5.646062435170393)
    ,(20, 1.779747248176307, 0.1599571284657202, 0.4619929909497245, 0.1350573392986281, 0.002505497200876161, 1.222686554736931, 0.0558012769147993, 0.470395379876706, 0.458505063791435, 1.782220711986362, 5.646135506337885)
    ,(21, 1.893116953282629, 0.1788852720546572, 0.463344889006791, 0.1342141940220571, 0.002501827854581418, 1.215554863078253, 0.05698702852100723, 0.4762830776872278, 0.4781467412798861, 1.76154656601274, 5.649800593850671)
    ,(22, 1.967607648748583, 0.1912188451219379, 0.4643864547853696, 0.1336919225237905, 0.002499522011528193, 1.211411067638892, 0.05775488007419582, 0.4802834613201696, 0.4911135027716114, 1.748574789597559, 5.654745450286327)
    ,(23, 1.927814677552533, 0.1846396157820891, 0.4638161154208973, 0.13396803995425, 0.002500744264033562, 1.213575352726281, 0.05734571444864756, 0.4781344718226137, 0.4841811050607063, 1.755448658194346, 5.651872951886958)
    ,(24, 1.777901081094114, 0.1596471408697937, 0.4619735677501637, 0.135071606711425, 0.002505558737912718, 1.22281184629913, 0.05578177669000084, 0.4703017304136381, 0.4581862403431174, 1.782567735076757, 5.646118631175312)
    ,(25, 1.617943931600166, 0.1325334753366008, 0.4606668351328551, 0.1363859233831501, 0.002511150679933397, 1.234993867542258, 0.05406440696852499, 0.4625111223802207, 0.4307120616926596, 1.814138050514401, 5.650851973848262)
)

$cols = @(2, 3, 4, 5, 7, 9, 10, 12, 13, 14, 15)

foreach ($entry in $newData) {
    $rowNum = $entry[0]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $colNum = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Cells.Item($rowNum, $colNum).Value = $value
    }
}
